# daily auto push: 2026-01-20 13:53 UTC
# Insert a new log row (2026/01/20, 火, 19, 201) into the date-ordered
# log table on row 678, pushing the existing rows 678-719 down to 679-720.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 678 (shifts rows 678:719 -> 679:720,
# and the sheet's dimension grows from D719 to D720 automatically).
$ws.Rows.Item(678).Insert()

# Column A holds date-like text ("yyyy/mm/dd"). Assigning it directly would
# make Excel auto-convert it into a date serial number, so we prefix it
# with a quote to force text entry, then reset the cell style back to
# Normal so no stray "quote prefix" style/number-format is left behind.
$ws.Range("A678").Value = "'2026/01/20"
$ws.Range("A678").Style = "Normal"

$ws.Range("B678").Value = "火"
$ws.Range("C678").Value = 19
$ws.Range("D678").Value = 201
